$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column A (Timestamp): shift dates by +4 days ---
# Rows 2-97: 12.06.2025 -> 16.06.2025 (serial 45820 -> 45824)
# Rows 98-193: 13.06.2025 -> 17.06.2025 (serial 45821 -> 45825)
for ($r = 2; $r -le 193; $r++) {
    $cur = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 1).Value = $cur + 4
}

# --- Update column E (Lookup): shift date prefix in text ---
for ($r = 2; $r -le 97; $r++) {
    $cur = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 5).Value = $cur.Replace("12.06.2025", "16.06.2025")
}
for ($r = 98; $r -le 193; $r++) {
    $cur = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 5).Value = $cur.Replace("13.06.2025", "17.06.2025")
}

# --- Update column B (Notified Production) and C (Actual Production) ---
$ws.Cells.Item(18, 2).Value = 3
$ws.Cells.Item(19, 2).Value = 3
$ws.Cells.Item(20, 2).Value = 15
$ws.Cells.Item(21, 2).Value = 26
$ws.Cells.Item(21, 3).Value = 7
$ws.Cells.Item(22, 2).Value = 114
$ws.Cells.Item(22, 3).Value = 27
$ws.Cells.Item(23, 2).Value = 126
$ws.Cells.Item(23, 3).Value = 55
$ws.Cells.Item(24, 2).Value = 166
$ws.Cells.Item(24, 3).Value = 91
$ws.Cells.Item(25, 2).Value = 189
$ws.Cells.Item(25, 3).Value = 130
$ws.Cells.Item(26, 2).Value = 469
$ws.Cells.Item(26, 3).Value = 217
$ws.Cells.Item(27, 2).Value = 500
$ws.Cells.Item(27, 3).Value = 296
$ws.Cells.Item(28, 2).Value = 533
$ws.Cells.Item(28, 3).Value = 372
$ws.Cells.Item(29, 2).Value = 574
$ws.Cells.Item(29, 3).Value = 482
$ws.Cells.Item(30, 2).Value = 1056
$ws.Cells.Item(30, 3).Value = 636
$ws.Cells.Item(31, 2).Value = 1125
$ws.Cells.Item(31, 3).Value = 754
$ws.Cells.Item(32, 2).Value = 1196
$ws.Cells.Item(32, 3).Value = 852
$ws.Cells.Item(33, 2).Value = 1242
$ws.Cells.Item(33, 3).Value = 943
$ws.Cells.Item(34, 2).Value = 1613
$ws.Cells.Item(34, 3).Value = 1078
$ws.Cells.Item(35, 2).Value = 1651
$ws.Cells.Item(35, 3).Value = 1193
$ws.Cells.Item(36, 2).Value = 1710
$ws.Cells.Item(36, 3).Value = 1278
$ws.Cells.Item(37, 2).Value = 1750
$ws.Cells.Item(37, 3).Value = 1332
$ws.Cells.Item(38, 2).Value = 2025
$ws.Cells.Item(38, 3).Value = 1421
$ws.Cells.Item(39, 2).Value = 2055
$ws.Cells.Item(39, 3).Value = 1521
$ws.Cells.Item(40, 2).Value = 2087
$ws.Cells.Item(40, 3).Value = 1579
$ws.Cells.Item(41, 2).Value = 2114
$ws.Cells.Item(41, 3).Value = 1618
$ws.Cells.Item(42, 2).Value = 2284
$ws.Cells.Item(42, 3).Value = 1688
$ws.Cells.Item(43, 2).Value = 2301
$ws.Cells.Item(43, 3).Value = 1722
$ws.Cells.Item(44, 2).Value = 2316
$ws.Cells.Item(44, 3).Value = 1754
$ws.Cells.Item(45, 2).Value = 2328
$ws.Cells.Item(45, 3).Value = 1783
$ws.Cells.Item(46, 2).Value = 2392
$ws.Cells.Item(46, 3).Value = 1807
$ws.Cells.Item(47, 2).Value = 2399
$ws.Cells.Item(47, 3).Value = 1827
$ws.Cells.Item(48, 2).Value = 2405
$ws.Cells.Item(48, 3).Value = 1843
$ws.Cells.Item(49, 2).Value = 2407
$ws.Cells.Item(49, 3).Value = 1837
$ws.Cells.Item(50, 2).Value = 2400
$ws.Cells.Item(50, 3).Value = 1845
$ws.Cells.Item(51, 2).Value = 2398
$ws.Cells.Item(51, 3).Value = 1842
$ws.Cells.Item(52, 2).Value = 2395
$ws.Cells.Item(52, 3).Value = 1825
$ws.Cells.Item(53, 2).Value = 2388
$ws.Cells.Item(53, 3).Value = 1809
$ws.Cells.Item(54, 2).Value = 2333
$ws.Cells.Item(54, 3).Value = 1781
$ws.Cells.Item(55, 2).Value = 2323
$ws.Cells.Item(55, 3).Value = 1775
$ws.Cells.Item(56, 2).Value = 2312
$ws.Cells.Item(56, 3).Value = 1754
$ws.Cells.Item(57, 2).Value = 2299
$ws.Cells.Item(57, 3).Value = 1721
$ws.Cells.Item(58, 2).Value = 2179
$ws.Cells.Item(58, 3).Value = 1649
$ws.Cells.Item(59, 2).Value = 2160
$ws.Cells.Item(59, 3).Value = 1622
$ws.Cells.Item(60, 2).Value = 2139
$ws.Cells.Item(60, 3).Value = 1610
$ws.Cells.Item(61, 2).Value = 2116
$ws.Cells.Item(61, 3).Value = 1582
$ws.Cells.Item(62, 2).Value = 1890
$ws.Cells.Item(62, 3).Value = 1512
$ws.Cells.Item(63, 2).Value = 1863
$ws.Cells.Item(63, 3).Value = 1480
$ws.Cells.Item(64, 2).Value = 1833
$ws.Cells.Item(64, 3).Value = 1377
$ws.Cells.Item(65, 2).Value = 1802
$ws.Cells.Item(65, 3).Value = 1286
$ws.Cells.Item(66, 2).Value = 1487
$ws.Cells.Item(66, 3).Value = 1158
$ws.Cells.Item(67, 2).Value = 1444
$ws.Cells.Item(67, 3).Value = 1123
$ws.Cells.Item(68, 2).Value = 1400
$ws.Cells.Item(68, 3).Value = 1074
$ws.Cells.Item(69, 2).Value = 1359
$ws.Cells.Item(69, 3).Value = 986
$ws.Cells.Item(70, 2).Value = 912
$ws.Cells.Item(70, 3).Value = 836
$ws.Cells.Item(71, 2).Value = 870
$ws.Cells.Item(71, 3).Value = 756
$ws.Cells.Item(72, 2).Value = 828
$ws.Cells.Item(72, 3).Value = 638
$ws.Cells.Item(73, 2).Value = 790
$ws.Cells.Item(73, 3).Value = 560
$ws.Cells.Item(74, 2).Value = 373
$ws.Cells.Item(74, 3).Value = 423
$ws.Cells.Item(75, 2).Value = 339
$ws.Cells.Item(75, 3).Value = 321
$ws.Cells.Item(76, 2).Value = 313
$ws.Cells.Item(76, 3).Value = 243
$ws.Cells.Item(77, 2).Value = 290
$ws.Cells.Item(77, 3).Value = 183
$ws.Cells.Item(78, 2).Value = 85
$ws.Cells.Item(78, 3).Value = 125
$ws.Cells.Item(79, 2).Value = 70
$ws.Cells.Item(79, 3).Value = 84
$ws.Cells.Item(80, 2).Value = 60
$ws.Cells.Item(80, 3).Value = 45
$ws.Cells.Item(81, 2).Value = 52
$ws.Cells.Item(81, 3).Value = 17
$ws.Cells.Item(82, 2).Value = 8
$ws.Cells.Item(82, 3).Value = 6
$ws.Cells.Item(83, 2).Value = 8
$ws.Cells.Item(83, 3).Value = 2
$ws.Cells.Item(84, 2).Value = 8
$ws.Cells.Item(85, 2).Value = 8
$ws.Cells.Item(86, 2).Value = 3
$ws.Cells.Item(90, 2).Value = 2
$ws.Cells.Item(114, 2).Value = 6
$ws.Cells.Item(115, 2).Value = 6
$ws.Cells.Item(116, 2).Value = 8
$ws.Cells.Item(117, 2).Value = 12
$ws.Cells.Item(117, 3).Value = 8
$ws.Cells.Item(118, 2).Value = 113
$ws.Cells.Item(118, 3).Value = 28
$ws.Cells.Item(119, 2).Value = 125
$ws.Cells.Item(119, 3).Value = 61
$ws.Cells.Item(120, 2).Value = 142
$ws.Cells.Item(120, 3).Value = 99
$ws.Cells.Item(121, 2).Value = 164
$ws.Cells.Item(121, 3).Value = 145
$ws.Cells.Item(122, 2).Value = 464
$ws.Cells.Item(122, 3).Value = 229
$ws.Cells.Item(123, 2).Value = 493
$ws.Cells.Item(123, 3).Value = 316
$ws.Cells.Item(124, 2).Value = 527
$ws.Cells.Item(124, 3).Value = 389
$ws.Cells.Item(125, 2).Value = 567
$ws.Cells.Item(125, 3).Value = 470
$ws.Cells.Item(126, 2).Value = 1071
$ws.Cells.Item(126, 3).Value = 599
$ws.Cells.Item(127, 2).Value = 1113
$ws.Cells.Item(127, 3).Value = 729
$ws.Cells.Item(128, 2).Value = 1160
$ws.Cells.Item(128, 3).Value = 828
$ws.Cells.Item(129, 2).Value = 1208
$ws.Cells.Item(129, 3).Value = 898
$ws.Cells.Item(130, 2).Value = 1589
$ws.Cells.Item(130, 3).Value = 1021
$ws.Cells.Item(131, 2).Value = 1631
$ws.Cells.Item(131, 3).Value = 1138
$ws.Cells.Item(132, 2).Value = 1673
$ws.Cells.Item(132, 3).Value = 1250
$ws.Cells.Item(133, 2).Value = 1714
$ws.Cells.Item(133, 3).Value = 1324
$ws.Cells.Item(134, 2).Value = 1959
$ws.Cells.Item(134, 3).Value = 1428
$ws.Cells.Item(135, 2).Value = 2004
$ws.Cells.Item(135, 3).Value = 1498
$ws.Cells.Item(136, 2).Value = 2036
$ws.Cells.Item(136, 3).Value = 1540
$ws.Cells.Item(137, 2).Value = 2062
$ws.Cells.Item(137, 3).Value = 1560
$ws.Cells.Item(138, 2).Value = 2199
$ws.Cells.Item(139, 2).Value = 2214
$ws.Cells.Item(140, 2).Value = 2227
$ws.Cells.Item(141, 2).Value = 2239
$ws.Cells.Item(142, 2).Value = 2293
$ws.Cells.Item(143, 2).Value = 2299
$ws.Cells.Item(144, 2).Value = 2304
$ws.Cells.Item(145, 2).Value = 2303
$ws.Cells.Item(146, 2).Value = 2288
$ws.Cells.Item(147, 2).Value = 2285
$ws.Cells.Item(148, 2).Value = 2280
$ws.Cells.Item(149, 2).Value = 2272
$ws.Cells.Item(150, 2).Value = 2160
$ws.Cells.Item(151, 2).Value = 2149
$ws.Cells.Item(152, 2).Value = 2136
$ws.Cells.Item(153, 2).Value = 2124
$ws.Cells.Item(154, 2).Value = 2004
$ws.Cells.Item(155, 2).Value = 1981
$ws.Cells.Item(156, 2).Value = 1958
$ws.Cells.Item(157, 2).Value = 1932
$ws.Cells.Item(158, 2).Value = 1708
$ws.Cells.Item(159, 2).Value = 1674
$ws.Cells.Item(160, 2).Value = 1643
$ws.Cells.Item(161, 2).Value = 1610
$ws.Cells.Item(162, 2).Value = 1318
$ws.Cells.Item(163, 2).Value = 1274
$ws.Cells.Item(164, 2).Value = 1231
$ws.Cells.Item(165, 2).Value = 1196
$ws.Cells.Item(166, 2).Value = 796
$ws.Cells.Item(167, 2).Value = 746
$ws.Cells.Item(168, 2).Value = 708
$ws.Cells.Item(169, 2).Value = 676
$ws.Cells.Item(170, 2).Value = 309
$ws.Cells.Item(171, 2).Value = 280
$ws.Cells.Item(172, 2).Value = 256
$ws.Cells.Item(173, 2).Value = 237
$ws.Cells.Item(174, 2).Value = 66
$ws.Cells.Item(175, 2).Value = 53
$ws.Cells.Item(176, 2).Value = 44
$ws.Cells.Item(177, 2).Value = 38
$ws.Cells.Item(179, 2).Value = 5
$ws.Cells.Item(180, 2).Value = 5
$ws.Cells.Item(181, 2).Value = 5
$ws.Cells.Item(186, 2).Value = 2
